$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Drive Team Data")

$ws.Range("J5").Value = -3.1963437456007147
$ws.Range("T5").Value = -3.5567052061310243
$ws.Range("J6").Value = 1.4450645858110902
$ws.Range("J7").Value = -0.089635456729177054
$ws.Range("T7").Value = 6.9088004812459207
$ws.Range("J8").Value = -0.44990728864085194
$ws.Range("T8").Value = -11.46303565018658
$ws.Range("J9").Value = 4.1915010427709571
$ws.Range("T9").Value = 5.2001634839811572
$ws.Range("J10").Value = 2.656801000230689
$ws.Range("J11").Value = -0.54862380879740957
$ws.Range("J12").Value = 4.092784522614398
$ws.Range("J13").Value = 2.558084480074128
$ws.Range("T13").Value = 2.2930925827348165
$ws.Range("J16").Value = -4.3750300198451972
$ws.Range("T16").Value = -20.664353527672162
$ws.Range("J17").Value = -0.66190335471574802
$ws.Range("J18").Value = -1.8896633887479708
$ws.Range("T18").Value = 6.9305881518404489
$ws.Range("J19").Value = -2.1778808542773049
$ws.Range("J20").Value = 1.5352458108521443
$ws.Range("J21").Value = 0.30748577681992728
$ws.Range("J22").Value = -2.256854070402551
$ws.Range("J23").Value = 1.4562725947268966
$ws.Range("J24").Value = 0.22851256069468029
$ws.Range("T24").Value = 5.0789414150025314
$ws.Range("J27").Value = -1.6172372164257025
$ws.Range("J28").Value = 2.0958894487037396
$ws.Range("J29").Value = 0.86812941467152915
$ws.Range("J30").Value = 0.57991194914218891
$ws.Range("J31").Value = 4.2930386142716319
$ws.Range("T31").Value = 5.0649884562557341
$ws.Range("J32").Value = 3.0652785802394238
$ws.Range("J33").Value = 0.50093873301694769
$ws.Range("J34").Value = 4.2140653981463911
$ws.Range("J35").Value = 2.9863053641141732
$ws.Range("T35").Value = -4.3566944184023662
$ws.Range("J38").Value = -1.8479996324941204
$ws.Range("T38").Value = 9.7555614361352863
$ws.Range("J39").Value = 1.865127032635322
$ws.Range("J40").Value = 0.63736699860310919
$ws.Range("J41").Value = 0.34914953307376878
$ws.Range("T41").Value = -11.46303565018658
$ws.Range("J42").Value = 4.0622761982032118
$ws.Range("J43").Value = 2.8345161641710037
$ws.Range("J44").Value = 0.27017631694852273
$ws.Range("J45").Value = 3.9833029820779693
$ws.Range("J46").Value = 2.7555429480457549
$ws.Range("T46").Value = -0.3807357396422919

$ws.Range("V5").Select()

$wsReturn = $wb.Worksheets.Item("Per Member Data")
$wsReturn.Activate()
